# Update TestRunner test data:
#  - C2 (Execution for test case 101) gains the Regression group too
#  - Test case 105 (row 6) moves from the Regression group into Sanity
#  - Two new test cases (110, 111) are appended to the Sanity group

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Test case 101's Execution column now also covers Regression
$ws.Range("C2").Value = "Groups=Smoke,Regression"

# Test case 105 switches from Regression to Sanity
$ws.Range("B6").Value = "Sanity"

# Append new test cases 110 and 111 to the Sanity group
$ws.Range("A11").Value = 110
$ws.Range("B11").Value = "Sanity"

$ws.Range("A12").Value = 111
$ws.Range("B12").Value = "Sanity"
